$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation with a leading apostrophe so Excel COM does not
# coerce numeric-looking strings (prices) into Double/Date values, matching
# the original inline-string (text) cell contents exactly.

$ws.Range('D2').Value = "'37.715.26"
$ws.Range('E2').Value = "'  -1.03%  "
$ws.Range('D3').Value = "'2.027.14"
$ws.Range('E3').Value = "'  -1.61%  "
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('D5').Value = "'227.01"
$ws.Range('E5').Value = "'  -1.52%  "
$ws.Range('E6').Value = "'  -0.49%  "
$ws.Range('D7').Value = "'59.18"
$ws.Range('E7').Value = "'  +1.40%  "
$ws.Range('E8').Value = "'  +0.06%  "
$ws.Range('E9').Value = "'  -1.06%  "
$ws.Range('D10').Value = "'0.0813"
$ws.Range('E10').Value = "'  +0.42%  "
$ws.Range('E11').Value = "'  -0.06%  "
$ws.Range('D12').Value = "'14.53"
$ws.Range('E12').Value = "'  -0.57%  "
$ws.Range('D13').Value = "'2.329.13"
$ws.Range('E13').Value = "'  -1.59%  "
$ws.Range('D14').Value = "'20.93"
$ws.Range('E14').Value = "'  +1.25%  "
$ws.Range('D15').Value = "'0.757"
$ws.Range('E15').Value = "'  +0.44%  "
$ws.Range('E16').Value = "'  -1.86%  "
$ws.Range('D17').Value = "'2.042.80"
$ws.Range('E17').Value = "'  -0.94%  "
$ws.Range('D18').Value = "'37.699.03"
$ws.Range('E18').Value = "'  -0.80%  "
$ws.Range('D19').Value = "'6.01"
$ws.Range('E19').Value = "'  -2.03%  "
$ws.Range('D20').Value = "'69.81"
$ws.Range('E20').Value = "'  -0.18%  "
$ws.Range('D21').Value = "'0.0₃0822"
$ws.Range('E21').Value = "'  -1.11%  "
$ws.Range('D22').Value = "'224.70"
$ws.Range('E22').Value = "'  -0.07%  "
$ws.Range('E23').Value = "'  -0.02%  "
$ws.Range('E24').Value = "'  -2.77%  "
$ws.Range('E25').Value = "'  -1.70%  "
$ws.Range('E26').Value = "'  -0.92%  "
$ws.Range('D27').Value = "'165.02"
$ws.Range('E27').Value = "'  -0.84%  "
$ws.Range('D28').Value = "'0.128"
$ws.Range('E28').Value = "'  -2.89%  "
$ws.Range('D29').Value = "'18.90"
$ws.Range('E29').Value = "'  -0.79%  "
$ws.Range('E30').Value = "'  -5.10%  "
$ws.Range('E31').Value = "'  +0.95%  "
$ws.Range('E32').Value = "'  -2.88%  "
$ws.Range('E33').Value = "'  +4.47%  "
$ws.Range('B34').Value = "'Hedera"
$ws.Range('C34').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D34').Value = "'0.0600"
$ws.Range('E34').Value = "'  -2.06%  "
$ws.Range('B35').Value = "'InternetComputer(DFINITY)"
$ws.Range('C35').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D35').Value = "'4.47"
$ws.Range('E35').Value = "'  -3.08%  "
$ws.Range('D36').Value = "'6.34"
$ws.Range('E36').Value = "'  +5.35%  "
$ws.Range('D37').Value = "'2.24"
$ws.Range('E37').Value = "'  -4.26%  "
$ws.Range('D38').Value = "'3.22"
$ws.Range('E38').Value = "'  -2.90%  "
$ws.Range('E39').Value = "'  +0.01%  "
$ws.Range('D40').Value = "'1.530.64"
$ws.Range('E40').Value = "'  +3.40%  "
$ws.Range('E41').Value = "'  -0.81%  "
$ws.Range('B42').Value = "'Aave"
$ws.Range('C42').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D42').Value = "'96.60"
$ws.Range('E42').Value = "'  -1.86%  "
$ws.Range('B43').Value = "'InjectiveProtocol"
$ws.Range('C43').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('D43').Value = "'16.72"
$ws.Range('E43').Value = "'  +0.19%  "
$ws.Range('B44').Value = "'HuobiToken"
$ws.Range('C44').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D44').Value = "'2.84"
$ws.Range('E44').Value = "'  -0.57%  "
$ws.Range('B45').Value = "'Cronos"
$ws.Range('C45').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D45').Value = "'0.0915"
$ws.Range('E45').Value = "'  -2.80%  "
$ws.Range('B46').Value = "'FTXToken"
$ws.Range('C46').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D46').Value = "'4.22"
$ws.Range('E46').Value = "'  +2.84%  "
$ws.Range('E47').Value = "'  -1.75%  "
$ws.Range('E48').Value = "'  -1.66%  "
$ws.Range('E49').Value = "'  -0.45%  "
$ws.Range('D50').Value = "'7.10"
$ws.Range('E50').Value = "'  +0.24%  "
$ws.Range('D51').Value = "'2.218.42"
$ws.Range('E51').Value = "'  -1.62%  "
